# Added 1.1.0 of term
# - Bump the "Version" property from 1.0.0 to 1.1.0
# - Bump the "Date" property to the new publication timestamp
# - Re-apply the wrap/vertical-top alignment on the bordered table cells so
#   the style definitions explicitly record applyAlignment="true"
#   (this mirrors the tool that produced the original workbook re-writing
#   the cell styles when the sheet was regenerated for the new version).

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"

# Header row + data rows of the Metadata sheet share the two bordered
# cell styles (bold header / plain body) that need applyAlignment="true".
$meta.Range("A1:B1").WrapText = $true
$meta.Range("A2:B14").WrapText = $true

# The second sheet reuses the very same two styles, so touch its cells too.
$fsiii = $wb.Worksheets.Item("Include from FSIII")
$fsiii.Range("A1:C1").WrapText = $true
$fsiii.Range("A2:B4").WrapText = $true
$fsiii.Range("C2").WrapText = $true
